$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "id" column (A) - data shifts left so "name" becomes column A
$ws.Columns.Item(1).Delete()

# Remove the duplicated date columns (now G:H after the shift) that held
# the "Jan 19, 2025" values - these are no longer exported
$ws.Range("G1:H3").Delete()

# Reflect the new selection left behind on the sheet (column G header)
$ws.Range("G1:G1048576").Select()
